# Fruta / hortaliza, semanal
# Update Fecha (D), Volumen (J), Precio minimo/maximo/promedio (K/L/M),
# Origen (O) and Precio $/Kg (P) for the Haba data rows (2-41).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44400
$ws.Range("K2").Value = 16500
$ws.Range("L2").Value = 16500
$ws.Range("M2").Value = 16500
$ws.Range("P2").Value = 660

$ws.Range("D3").Value = 44162
$ws.Range("J3").Value = 140
$ws.Range("K3").Value = 13000
$ws.Range("L3").Value = 13000
$ws.Range("M3").Value = 13000
$ws.Range("O3").Value = "Región del Maule"
$ws.Range("P3").Value = 520

$ws.Range("D4").Value = 44383
$ws.Range("K4").Value = 17000
$ws.Range("L4").Value = 17000
$ws.Range("M4").Value = 17000
$ws.Range("P4").Value = 680

$ws.Range("D5").Value = 44442
$ws.Range("J5").Value = 80
$ws.Range("K5").Value = 18000
$ws.Range("L5").Value = 18000
$ws.Range("M5").Value = 18000
$ws.Range("P5").Value = 720

$ws.Range("D7").Value = 44160
$ws.Range("J7").Value = 40
$ws.Range("K7").Value = 11500
$ws.Range("L7").Value = 11500
$ws.Range("M7").Value = 11500
$ws.Range("O7").Value = "Región del Maule"
$ws.Range("P7").Value = 460

$ws.Range("D8").Value = 44460
$ws.Range("J8").Value = 80
$ws.Range("K8").Value = 15000
$ws.Range("L8").Value = 15000
$ws.Range("M8").Value = 15000
$ws.Range("P8").Value = 600

$ws.Range("D9").Value = 44167
$ws.Range("J9").Value = 30
$ws.Range("K9").Value = 12000
$ws.Range("L9").Value = 12000
$ws.Range("M9").Value = 12000
$ws.Range("O9").Value = "Región de La Araucanía"
$ws.Range("P9").Value = 480

$ws.Range("D10").Value = 44365
$ws.Range("J10").Value = 80
$ws.Range("K10").Value = 20000
$ws.Range("L10").Value = 20000
$ws.Range("M10").Value = 20000
$ws.Range("P10").Value = 800

$ws.Range("D11").Value = 44376
$ws.Range("J11").Value = 70
$ws.Range("K11").Value = 17000
$ws.Range("L11").Value = 17000
$ws.Range("M11").Value = 17000
$ws.Range("O11").Value = "Provincia de Limarí"
$ws.Range("P11").Value = 680

$ws.Range("D12").Value = 44446
$ws.Range("J12").Value = 90

$ws.Range("D13").Value = 44379
$ws.Range("J13").Value = 70

$ws.Range("D14").Value = 44418
$ws.Range("J14").Value = 90
$ws.Range("K14").Value = 18000
$ws.Range("L14").Value = 18000
$ws.Range("M14").Value = 18000
$ws.Range("P14").Value = 720

$ws.Range("D15").Value = 44174
$ws.Range("J15").Value = 20
$ws.Range("K15").Value = 12000
$ws.Range("L15").Value = 12000
$ws.Range("M15").Value = 12000
$ws.Range("O15").Value = "Región de La Araucanía"
$ws.Range("P15").Value = 480

$ws.Range("D16").Value = 44449
$ws.Range("K16").Value = 17000
$ws.Range("L16").Value = 17000
$ws.Range("M16").Value = 17000
$ws.Range("P16").Value = 680

$ws.Range("D17").Value = 44463
$ws.Range("K17").Value = 16000
$ws.Range("L17").Value = 16000
$ws.Range("M17").Value = 16000
$ws.Range("P17").Value = 640

$ws.Range("D18").Value = 44428
$ws.Range("J18").Value = 80
$ws.Range("K18").Value = 18500
$ws.Range("L18").Value = 18500
$ws.Range("M18").Value = 18500
$ws.Range("P18").Value = 740

$ws.Range("D19").Value = 44421
$ws.Range("K19").Value = 17000
$ws.Range("L19").Value = 17000
$ws.Range("M19").Value = 17000
$ws.Range("P19").Value = 680

$ws.Range("D20").Value = 44169
$ws.Range("J20").Value = 160
$ws.Range("K20").Value = 13000
$ws.Range("L20").Value = 14000
$ws.Range("M20").Value = 13500
$ws.Range("P20").Value = 540

$ws.Range("D21").Value = 44411
$ws.Range("J21").Value = 80
$ws.Range("K21").Value = 18000
$ws.Range("L21").Value = 18000
$ws.Range("M21").Value = 18000
$ws.Range("O21").Value = "Provincia de Limarí"
$ws.Range("P21").Value = 720

$ws.Range("D22").Value = 44215
$ws.Range("J22").Value = 60
$ws.Range("K22").Value = 25000
$ws.Range("L22").Value = 25000
$ws.Range("M22").Value = 25000
$ws.Range("O22").Value = "Región Metropolitana"
$ws.Range("P22").Value = 1000

$ws.Range("D23").Value = 44203
$ws.Range("J23").Value = 50
$ws.Range("K23").Value = 22000
$ws.Range("L23").Value = 22000
$ws.Range("M23").Value = 22000
$ws.Range("O23").Value = "Región de La Araucanía"
$ws.Range("P23").Value = 880

$ws.Range("D24").Value = 44161
$ws.Range("J24").Value = 50
$ws.Range("K24").Value = 11500
$ws.Range("L24").Value = 11500
$ws.Range("M24").Value = 11500
$ws.Range("O24").Value = "Región del Maule"
$ws.Range("P24").Value = 460

$ws.Range("D25").Value = 44351
$ws.Range("J25").Value = 30
$ws.Range("K25").Value = 20000
$ws.Range("L25").Value = 20000
$ws.Range("M25").Value = 20000
$ws.Range("O25").Value = "Región Metropolitana"
$ws.Range("P25").Value = 800

$ws.Range("D26").Value = 44214
$ws.Range("J26").Value = 40
$ws.Range("K26").Value = 25000
$ws.Range("L26").Value = 25000
$ws.Range("M26").Value = 25000
$ws.Range("O26").Value = "Región de La Araucanía"
$ws.Range("P26").Value = 1000

$ws.Range("D27").Value = 44407
$ws.Range("K27").Value = 18000
$ws.Range("L27").Value = 18000
$ws.Range("M27").Value = 18000
$ws.Range("P27").Value = 720

$ws.Range("D28").Value = 44358
$ws.Range("J28").Value = 60
$ws.Range("K28").Value = 20000
$ws.Range("L28").Value = 20000
$ws.Range("M28").Value = 20000
$ws.Range("O28").Value = "Región Metropolitana"
$ws.Range("P28").Value = 800

$ws.Range("D29").Value = 44474
$ws.Range("J29").Value = 150

$ws.Range("D30").Value = 44362
$ws.Range("J30").Value = 60
$ws.Range("K30").Value = 20000
$ws.Range("L30").Value = 20000
$ws.Range("M30").Value = 20000
$ws.Range("P30").Value = 800

$ws.Range("D31").Value = 44435
$ws.Range("J31").Value = 170
$ws.Range("K31").Value = 18000
$ws.Range("L31").Value = 19000
$ws.Range("M31").Value = 18529
$ws.Range("O31").Value = "Provincia de Limarí"
$ws.Range("P31").Value = 741

$ws.Range("D32").Value = 44159
$ws.Range("J32").Value = 150

$ws.Range("D33").Value = 44166
$ws.Range("J33").Value = 120
$ws.Range("K33").Value = 12000
$ws.Range("L33").Value = 12000
$ws.Range("M33").Value = 12000
$ws.Range("P33").Value = 480

$ws.Range("D34").Value = 44386
$ws.Range("K34").Value = 17000
$ws.Range("L34").Value = 17000
$ws.Range("M34").Value = 17000
$ws.Range("P34").Value = 680

$ws.Range("D35").Value = 44369
$ws.Range("J35").Value = 70
$ws.Range("K35").Value = 18000
$ws.Range("L35").Value = 18000
$ws.Range("M35").Value = 18000
$ws.Range("O35").Value = "Provincia de Limarí"
$ws.Range("P35").Value = 720

$ws.Range("D36").Value = 44172
$ws.Range("J36").Value = 40
$ws.Range("K36").Value = 12000
$ws.Range("L36").Value = 12000
$ws.Range("M36").Value = 12000
$ws.Range("P36").Value = 480

$ws.Range("D37").Value = 44414
$ws.Range("J37").Value = 80
$ws.Range("K37").Value = 18000
$ws.Range("L37").Value = 18000
$ws.Range("M37").Value = 18000
$ws.Range("O37").Value = "Provincia de Limarí"
$ws.Range("P37").Value = 720

$ws.Range("D38").Value = 44392
$ws.Range("J38").Value = 90
$ws.Range("K38").Value = 16000
$ws.Range("L38").Value = 16000
$ws.Range("M38").Value = 16000
$ws.Range("P38").Value = 640

$ws.Range("D39").Value = 44425
$ws.Range("J39").Value = 90
$ws.Range("K39").Value = 18000
$ws.Range("L39").Value = 18000
$ws.Range("M39").Value = 18000
$ws.Range("P39").Value = 720

$ws.Range("D40").Value = 44390
$ws.Range("J40").Value = 80
$ws.Range("K40").Value = 16000
$ws.Range("L40").Value = 16000
$ws.Range("M40").Value = 16000
$ws.Range("O40").Value = "Provincia de Limarí"
$ws.Range("P40").Value = 640

$ws.Range("D41").Value = 44432
$ws.Range("J41").Value = 80
